# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker table (rows 16-21) was re-sorted: instead of being grouped by
# worker (each worker's two "Periodo Mora" rows 1603/1602 together), the
# data is now grouped by period (all 1602 rows first, then all 1603 rows).
# Only columns C (N Doc Trabajador), D (Nombre Trabajador) and E (Periodo
# Mora) change; B, F, G, H, I, J stay the same for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "45762493"
$ws.Range("D16").Value = "BEATRIZ ELENA PICO GUEVARA"
$ws.Range("E16").Value = "1602"

$ws.Range("C17").Value = "1067885958"
$ws.Range("D17").Value = "MARIA ANGELICA PIEDRAHITA QUINTANA"
$ws.Range("E17").Value = "1602"

$ws.Range("C18").Value = "1047480942"
$ws.Range("D18").Value = "LIETH EUGENIA GUERRERO CARDENAS"
$ws.Range("E18").Value = "1602"

$ws.Range("C19").Value = "45762493"
$ws.Range("D19").Value = "BEATRIZ ELENA PICO GUEVARA"
$ws.Range("E19").Value = "1603"

$ws.Range("C20").Value = "1067885958"
$ws.Range("D20").Value = "MARIA ANGELICA PIEDRAHITA QUINTANA"
$ws.Range("E20").Value = "1603"

$ws.Range("C21").Value = "1047480942"
$ws.Range("D21").Value = "LIETH EUGENIA GUERRERO CARDENAS"
$ws.Range("E21").Value = "1603"
